$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value2 = "Volume 31   Number  51"
$ws.Range("C9").Value2 = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- Data table updates (rows 14-33), refreshed weekly crime stats ---
# Style donor cells (unaffected by this edit) used to set correct cell type/format
# when a cell switches between numeric and "n/a" shared-string representation:
#   A39 -> style for text/"n/a" cells (s=13)
#   C39 -> style for integer-count cells (s=14)
#   K39 -> style for percent-change cells (s=15)
$ws.Range("N14").Value2 = -40
$ws.Range("A39").Copy($ws.Range("G15"))
$ws.Range("A39").Copy($ws.Range("H15"))
$ws.Range("A39").Copy($ws.Range("D16"))
$ws.Range("A39").Copy($ws.Range("E16"))
$ws.Range("F16").Value2 = 8
$ws.Range("G16").Value2 = 5
$ws.Range("H16").Value2 = 60
$ws.Range("I16").Value2 = 130
$ws.Range("K16").Value2 = -4.411764705882
$ws.Range("L16").Value2 = 8.333333333333
$ws.Range("M16").Value2 = -23.976608187134
$ws.Range("N16").Value2 = -82.456140350877
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = -66.666666666666
$ws.Range("F17").Value2 = 14
$ws.Range("H17").Value2 = -6.666666666666
$ws.Range("I17").Value2 = 251
$ws.Range("J17").Value2 = 222
$ws.Range("K17").Value2 = 13.063063063063
$ws.Range("L17").Value2 = 40.22346368715
$ws.Range("M17").Value2 = 87.31343283582
$ws.Range("N17").Value2 = -16.053511705685
$ws.Range("C39").Copy($ws.Range("C18"))
$ws.Range("C18").Value2 = 2
$ws.Range("C39").Copy($ws.Range("D18"))
$ws.Range("D18").Value2 = 3
$ws.Range("K39").Copy($ws.Range("E18"))
$ws.Range("E18").Value2 = -33.333333333333
$ws.Range("G18").Value2 = 6
$ws.Range("H18").Value2 = -33.333333333333
$ws.Range("I18").Value2 = 126
$ws.Range("J18").Value2 = 206
$ws.Range("K18").Value2 = -38.83495145631
$ws.Range("L18").Value2 = -31.521739130434
$ws.Range("M18").Value2 = -57.142857142857
$ws.Range("N18").Value2 = -92.288861689106
$ws.Range("C19").Value2 = 9
$ws.Range("D19").Value2 = 14
$ws.Range("E19").Value2 = -35.714285714285
$ws.Range("F19").Value2 = 33
$ws.Range("G19").Value2 = 51
$ws.Range("H19").Value2 = -35.294117647058
$ws.Range("I19").Value2 = 555
$ws.Range("J19").Value2 = 641
$ws.Range("K19").Value2 = -13.416536661466
$ws.Range("L19").Value2 = -22.594142259414
$ws.Range("M19").Value2 = 21.978021978022
$ws.Range("N19").Value2 = -27.450980392156
$ws.Range("C20").Value2 = 7
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = 133.333333333333
$ws.Range("F20").Value2 = 27
$ws.Range("G20").Value2 = 13
$ws.Range("H20").Value2 = 107.692307692308
$ws.Range("I20").Value2 = 220
$ws.Range("J20").Value2 = 160
$ws.Range("K20").Value2 = 37.5
$ws.Range("L20").Value2 = 50.684931506849
$ws.Range("M20").Value2 = 42.857142857142
$ws.Range("N20").Value2 = -87.757373400111
$ws.Range("D21").Value2 = 26
$ws.Range("E21").Value2 = -11.538461538461
$ws.Range("F21").Value2 = 86
$ws.Range("G21").Value2 = 91
$ws.Range("H21").Value2 = -5.494505494505
$ws.Range("I21").Value2 = 1304
$ws.Range("J21").Value2 = 1387
$ws.Range("K21").Value2 = -5.984138428262
$ws.Range("L21").Value2 = -4.678362573099
$ws.Range("M21").Value2 = 6.362153344208
$ws.Range("N21").Value2 = -75.256166982922
$ws.Range("A39").Copy($ws.Range("D22"))
$ws.Range("A39").Copy($ws.Range("E22"))
$ws.Range("M22").Value2 = -20.689655172413
$ws.Range("C24").Value2 = 30
$ws.Range("D24").Value2 = 26
$ws.Range("E24").Value2 = 15.384615384615
$ws.Range("F24").Value2 = 145
$ws.Range("G24").Value2 = 123
$ws.Range("H24").Value2 = 17.886178861788
$ws.Range("I24").Value2 = 1720
$ws.Range("J24").Value2 = 1651
$ws.Range("K24").Value2 = 4.179285281647
$ws.Range("L24").Value2 = -6.926406926406
$ws.Range("M24").Value2 = 79.166666666666
$ws.Range("C25").Value2 = 16
$ws.Range("D25").Value2 = 18
$ws.Range("E25").Value2 = -11.111111111111
$ws.Range("F25").Value2 = 87
$ws.Range("G25").Value2 = 76
$ws.Range("H25").Value2 = 14.473684210526
$ws.Range("I25").Value2 = 1019
$ws.Range("J25").Value2 = 965
$ws.Range("K25").Value2 = 5.595854922279
$ws.Range("L25").Value2 = -3.503787878787
$ws.Range("C26").Value2 = 8
$ws.Range("D26").Value2 = 14
$ws.Range("E26").Value2 = -42.857142857142
$ws.Range("F26").Value2 = 34
$ws.Range("G26").Value2 = 49
$ws.Range("H26").Value2 = -30.612244897959
$ws.Range("I26").Value2 = 601
$ws.Range("J26").Value2 = 616
$ws.Range("K26").Value2 = -2.435064935064
$ws.Range("L26").Value2 = 33.555555555555
$ws.Range("M26").Value2 = 41.411764705882
$ws.Range("A39").Copy($ws.Range("G27"))
$ws.Range("A39").Copy($ws.Range("H27"))
$ws.Range("A39").Copy($ws.Range("C28"))
$ws.Range("D28").Value2 = 1
$ws.Range("E28").Value2 = -100
$ws.Range("F28").Value2 = 1
$ws.Range("H28").Value2 = -66.666666666666
$ws.Range("J28").Value2 = 49
$ws.Range("K28").Value2 = -12.244897959183
$ws.Range("L28").Value2 = -15.686274509803
$ws.Range("M29").Value2 = -33.333333333333
$ws.Range("M30").Value2 = -66.666666666666
$ws.Range("C39").Copy($ws.Range("D31"))
$ws.Range("D31").Value2 = 2
$ws.Range("K39").Copy($ws.Range("E31"))
$ws.Range("E31").Value2 = -100
$ws.Range("A39").Copy($ws.Range("F31"))
$ws.Range("H31").Value2 = -100
$ws.Range("J31").Value2 = 8
$ws.Range("K31").Value2 = 100
$ws.Range("A39").Copy($ws.Range("G33"))
$ws.Range("A39").Copy($ws.Range("H33"))

# --- Insert blank spacer row so footer rows shift down (56->57, 57->58) ---
$ws.Rows.Item(56).Insert()
$ws.Rows.Item(56).Clear()
